$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 335; this shifts the existing rows 335-380 down to 336-381,
# carrying their data/formatting with them (matches the diff, which shows every
# row from 335 onward effectively taking on the values of the row above it, with
# a brand-new row of data appearing at 335 and the dimension growing to R381).
$ws.Rows.Item(335).Insert()

# Populate the newly inserted row 335 with the new record's data.
$ws.Range("A335").Value = 3
$ws.Range("B335").Value = "Femacal de La Calera"
$ws.Range("C335").Value = "Coquimbo"
$ws.Range("D335").Value = 44776
$ws.Range("E335").Value = 5
$ws.Range("F335").Value = 100112012
$ws.Range("G335").Value = "Espinaca"
$ws.Range("H335").Value = "Sin especificar"
$ws.Range("I335").Value = "Primera"
$ws.Range("J335").Value = 170
$ws.Range("K335").Value = 4000
$ws.Range("L335").Value = 4500
$ws.Range("M335").Value = 4265
$ws.Range("N335").Value = "$/docena de atados (3 kilos)"
$ws.Range("O335").Value = "Provincia de Quillota"
$ws.Range("P335").Value = 1422
$ws.Range("Q335").Value = 3
$ws.Range("R335").Value = "Hortaliza"
